$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18: fill in column A (the "file name" cell for the new block)
$ws.Range("A18").Value = "SCRIPT/T01P01A/um2508.ssb"

# New dialogue block spanning rows 19-20 (written in shared-string creation order)
$ws.Range("C19").Value = ' You\''ve brought peace to the\nworld at last!'
$ws.Range("C20").Value = ' Thank you!'

$ws.Range("A19").Value = "SCRIPT/T01P02A/us0105.ssb"
$ws.Range("D19").Value = ' Наконец-то вы смогли принести\nпокой в этот мир!'
$ws.Range("D20").Value = ' Спасибо вам!'

$ws.Range("E19").Value = ' Îàëïîåø-óï âú òíïãìé ðñéîåòóé\nðïëïê â üóïó íéñ!'
$ws.Range("E20").Value = ' Òðàòéáï âàí!'

$ws.Range("B19").Value = 66
$ws.Range("B20").Value = 69

# Row heights matching the new content (wrapped, taller rows)
$ws.Rows.Item(18).RowHeight = 43.2
$ws.Rows.Item(19).RowHeight = 43.2

# Update view: scroll down and select E20, matching the saved view state
$excel.ActiveWindow.ScrollRow = 16
$ws.Range("E20").Select() | Out-Null
